$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B1").Value = "pess"
$ws.Range("C1").Value = "X0"
$ws.Range("D1").Value = "X20"
$ws.Range("E1").Value = "X40"
$ws.Range("F1").Value = "X60"
$ws.Range("G1").Value = "X80"
$ws.Range("H1").Value = "X100"
$ws.Range("I1").Value = "opt"
